# Auto-generated edit script: update cryptos price/volume table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "51.287.95"
Set-TextValue $ws.Range("E2") "  -0.59%  "

Set-TextValue $ws.Range("D3") "2.973.18"
Set-TextValue $ws.Range("E3") "  +0.19%  "

Set-TextValue $ws.Range("D4") "1.00"
Set-TextValue $ws.Range("E4") "  +0.17%  "

Set-TextValue $ws.Range("D5") "381.89"
Set-TextValue $ws.Range("E5") "  +1.15%  "

Set-TextValue $ws.Range("D6") "102.08"
Set-TextValue $ws.Range("E6") "  -2.82%  "

Set-TextValue $ws.Range("D7") "0.541"
Set-TextValue $ws.Range("E7") "  -1.04%  "

Set-TextValue $ws.Range("E8") "  +0.16%  "

Set-TextValue $ws.Range("D9") "0.589"
Set-TextValue $ws.Range("E9") "  -1.48%  "

Set-TextValue $ws.Range("D10") "36.83"
Set-TextValue $ws.Range("E10") "  -1.96%  "

Set-TextValue $ws.Range("E11") "  -0.40%  "

Set-TextValue $ws.Range("D12") "0.0840"
Set-TextValue $ws.Range("E12") "  -0.52%  "

Set-TextValue $ws.Range("D13") "3.443.40"
Set-TextValue $ws.Range("E13") "  +0.62%  "

Set-TextValue $ws.Range("D14") "18.13"
Set-TextValue $ws.Range("E14") "  -2.20%  "

Set-TextValue $ws.Range("E15") "  +0.70%  "

Set-TextValue $ws.Range("D16") "2.977.36"
Set-TextValue $ws.Range("E16") "  +1.03%  "

Set-TextValue $ws.Range("D17") "0.990"
Set-TextValue $ws.Range("E17") "  +4.42%  "

Set-TextValue $ws.Range("D18") "51.225.68"
Set-TextValue $ws.Range("E18") "  -0.76%  "

Set-TextValue $ws.Range("D19") "3.24"
Set-TextValue $ws.Range("E19") "  -2.97%  "

Set-TextValue $ws.Range("D20") "7.36"
Set-TextValue $ws.Range("E20") "  -0.31%  "

Set-TextValue $ws.Range("D21") "12.77"
Set-TextValue $ws.Range("E21") "  -2.89%  "

Set-TextValue $ws.Range("D22") "0.0₃0955"
Set-TextValue $ws.Range("E22") "  +0.16%  "

Set-TextValue $ws.Range("D23") "68.81"
Set-TextValue $ws.Range("E23") "  -0.12%  "

Set-TextValue $ws.Range("D24") "260.50"
Set-TextValue $ws.Range("E24") "  -0.90%  "

Set-TextValue $ws.Range("D25") "2.89"
Set-TextValue $ws.Range("E25") "  +6.12%  "

Set-TextValue $ws.Range("D26") "8.08"
Set-TextValue $ws.Range("E26") "  +12.42%  "

Set-TextValue $ws.Range("D27") "7.41"
Set-TextValue $ws.Range("E27") "  +6.21%  "

Set-TextValue $ws.Range("D28") "0.117"
Set-TextValue $ws.Range("E28") "  +12.13%  "

Set-TextValue $ws.Range("E29") "  -0.12%  "

Set-TextValue $ws.Range("E30") "  -2.39%  "

Set-TextValue $ws.Range("E31") "  -0.03%  "

Set-TextValue $ws.Range("D32") "25.84"
Set-TextValue $ws.Range("E32") "  -0.86%  "

Set-TextValue $ws.Range("D33") "9.81"
Set-TextValue $ws.Range("E33") "  -1.51%  "

Set-TextValue $ws.Range("D34") "34.27"
Set-TextValue $ws.Range("E34") "  -2.24%  "

Set-TextValue $ws.Range("D35") "50.82"
Set-TextValue $ws.Range("E35") "  +1.04%  "

Set-TextValue $ws.Range("E36") "  -2.58%  "

Set-TextValue $ws.Range("D37") "0.0452"
Set-TextValue $ws.Range("E37") "  +4.55%  "

Set-TextValue $ws.Range("E38") "  -0.09%  "

Set-TextValue $ws.Range("D39") "2.96"
Set-TextValue $ws.Range("E39") "  -3.19%  "

Set-TextValue $ws.Range("D40") "16.96"
Set-TextValue $ws.Range("E40") "  -1.95%  "

Set-TextValue $ws.Range("D41") "2.55"
Set-TextValue $ws.Range("E41") "  -3.05%  "

Set-TextValue $ws.Range("E42") "  +0.77%  "

Set-TextValue $ws.Range("D43") "1.81"
Set-TextValue $ws.Range("E43") "  -3.20%  "

Set-TextValue $ws.Range("D44") "122.75"
Set-TextValue $ws.Range("E44") "  +2.18%  "

Set-TextValue $ws.Range("D45") "21.39"
Set-TextValue $ws.Range("E45") "  -3.88%  "

Set-TextValue $ws.Range("E46") "  -1.14%  "

Set-TextValue $ws.Range("D47") "0.271"
Set-TextValue $ws.Range("E47") "  +0.81%  "

# Rows 49/50: coins swapped (NEARProtocol <-> Maker) with refreshed price/volume data
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue $ws.Range("D49") "2.026.70"
Set-TextValue $ws.Range("E49") "  -0.94%  "

$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D50") "3.26"
Set-TextValue $ws.Range("E50") "  +0.92%  "

Set-TextValue $ws.Range("D51") "0.0331"
Set-TextValue $ws.Range("E51") "  +0.67%  "
